$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 437, shifting the existing row 437 (and everything
# below it) down by one. This grows the used range from R526 to R527.
$ws.Rows("437:437").Insert()

# Populate the newly inserted row 437 with this week's price observation
# for Perejil at Vega Central Mapocho de Santiago.
$ws.Cells.Item(437, 1).Value  = 9
$ws.Cells.Item(437, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(437, 3).Value  = "Metropolitana"
$ws.Cells.Item(437, 4).Value  = 44995
$ws.Cells.Item(437, 5).Value  = 13
$ws.Cells.Item(437, 6).Value  = 100112044
$ws.Cells.Item(437, 7).Value  = "Perejil"
$ws.Cells.Item(437, 8).Value  = "Sin especificar"
$ws.Cells.Item(437, 9).Value  = "Primera"
$ws.Cells.Item(437, 10).Value = 70
$ws.Cells.Item(437, 11).Value = 13000
$ws.Cells.Item(437, 12).Value = 14000
$ws.Cells.Item(437, 13).Value = 13500
$ws.Cells.Item(437, 14).Value = "$/docena de atados"
$ws.Cells.Item(437, 15).Value = "Región Metropolitana"
$ws.Cells.Item(437, 16).Value = 4500
$ws.Cells.Item(437, 17).Value = 3
$ws.Cells.Item(437, 18).Value = "Hortaliza"
